# Chapter 6 - Example 5: batch max/min statistics
# Adds "最大销售利润" (max sales profit) and "最小销售利润" (min sales profit)
# summary cells next to the existing 产品销售统计表 (product sales stats) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary labels + MAX/MIN values of the "销售利润" column (H2:H13) ---
$ws.Range("I1").Value = "最大销售利润"
$ws.Range("J1").Value = $ws.Application.WorksheetFunction.Max($ws.Range("H2:H13"))

$ws.Range("I2").Value = "最小销售利润"
$ws.Range("J2").Value = $ws.Application.WorksheetFunction.Min($ws.Range("H2:H13"))

# Give J1/J2 the same currency number format already used by column H
# (copy/paste-format reuses the existing style record instead of minting a new one)
$ws.Range("H2").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column widths: mirror the "best fit" autosize Excel applied after adding I:J ---
$ws.Columns.Item(1).ColumnWidth = 7.857142857142857
$ws.Columns.Item(2).ColumnWidth = 7.571428571428571
$ws.Columns.Item(3).ColumnWidth = 13.571428571428571
$ws.Columns.Item(4).ColumnWidth = 13.571428571428571
$ws.Columns.Item(5).ColumnWidth = 12.857142857142858
$ws.Columns.Item(6).ColumnWidth = 9.714285714285714
$ws.Columns.Item(7).ColumnWidth = 9.714285714285714
$ws.Columns.Item(8).ColumnWidth = 9.714285714285714
$ws.Columns.Item(9).ColumnWidth = 11.0
$ws.Columns.Item(10).ColumnWidth = 9.714285714285714
